$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit adds averaged-intensity results for three new spiral sampling
# schemes. Gaussian-Quadrature (previously the last data row) now comes right
# after the "Ring Perpendicular to TD" rows, followed by the three new Spiral
# rows, and the remaining schemes (NoRotation-tilt60deg .. HexGrid-60degTilt5degRes)
# shift down to make room.

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9834174039073463
$ws.Range("D10").Value = 0.9924703951312377
$ws.Range("E10").Value = 1.026900555561042
$ws.Range("F10").Value = 0.9834174039073463
$ws.Range("G10").Value = 0.9175800386854445
$ws.Range("H10").Value = 1.201031381252061
$ws.Range("I10").Value = 1.007901537378366
$ws.Range("J10").Value = 0.9924703951312377
$ws.Range("K10").Value = 1.00968547534614
$ws.Range("L10").Value = 0.996551439626743
$ws.Range("M10").Value = 1.021550218652583

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9401723248035061
$ws.Range("D11").Value = 1.889623144686522
$ws.Range("E11").Value = 0.7018675165841424
$ws.Range("F11").Value = 0.9401723248035061
$ws.Range("G11").Value = 1.4660457232924
$ws.Range("H11").Value = 0.328386906376964
$ws.Range("I11").Value = 0.7650490148542836
$ws.Range("J11").Value = 1.889623144686522
$ws.Range("K11").Value = 1.295745330635332
$ws.Range("L11").Value = 1.117958827719419
$ws.Range("M11").Value = 1.015190771766303

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9374496479028154
$ws.Range("D12").Value = 1.893761146967791
$ws.Range("E12").Value = 0.7015036444599082
$ws.Range("F12").Value = 0.9374496479028154
$ws.Range("G12").Value = 1.46821720311783
$ws.Range("H12").Value = 0.3274990554223933
$ws.Range("I12").Value = 0.764229370257537
$ws.Range("J12").Value = 1.893761146967791
$ws.Range("K12").Value = 1.297632395713849
$ws.Range("L12").Value = 1.117541021808333
$ws.Range("M12").Value = 1.015443344688046

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9396232570599472
$ws.Range("D13").Value = 1.890164067852492
$ws.Range("E13").Value = 0.7018487639650598
$ws.Range("F13").Value = 0.9396232570599472
$ws.Range("G13").Value = 1.466471032416336
$ws.Range("H13").Value = 0.3282244278643351
$ws.Range("I13").Value = 0.7649764300113838
$ws.Range("J13").Value = 1.890164067852492
$ws.Range("K13").Value = 1.296006415908776
$ws.Range("L13").Value = 1.117814836484361
$ws.Range("M13").Value = 1.015217996528259

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.408632
$ws.Range("D14").Value = 1.190688000000001
$ws.Range("E14").Value = 0.7513119999999998
$ws.Range("F14").Value = 1.408632
$ws.Range("G14").Value = 1.1114
$ws.Range("H14").Value = 0.4694839999999998
$ws.Range("I14").Value = 0.8944919999999995
$ws.Range("J14").Value = 1.190688000000001
$ws.Range("K14").Value = 0.9710000000000003
$ws.Range("L14").Value = 1.189816
$ws.Range("M14").Value = 0.9710013333333333

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.97
$ws.Range("D15").Value = 0.21
$ws.Range("E15").Value = 0.8426125000000002
$ws.Range("F15").Value = 1.97
$ws.Range("G15").Value = 0.64
$ws.Range("H15").Value = 0.66
$ws.Range("I15").Value = 1.093187500000002
$ws.Range("J15").Value = 0.21
$ws.Range("K15").Value = 0.5263062500000001
$ws.Range("L15").Value = 1.248153125
$ws.Range("M15").Value = 0.9026333333333336

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.5545899380736
$ws.Range("D16").Value = 0.5397915715584015
$ws.Range("E16").Value = 0.9099897923583969
$ws.Range("F16").Value = 1.5545899380736
$ws.Range("G16").Value = 0.790208806297601
$ws.Range("H16").Value = 0.8104218732544011
$ws.Range("I16").Value = 1.054980195840001
$ws.Range("J16").Value = 0.5397915715584015
$ws.Range("K16").Value = 0.7248906819583992
$ws.Range("L16").Value = 1.139740310016
$ws.Range("M16").Value = 0.943330362897067

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9954199582472592
$ws.Range("D17").Value = 0.9933137835985731
$ws.Range("E17").Value = 0.9942538958394512
$ws.Range("F17").Value = 0.9954199582472592
$ws.Range("G17").Value = 0.9934133634755241
$ws.Range("H17").Value = 0.9936876058048703
$ws.Range("I17").Value = 0.9956971022920688
$ws.Range("J17").Value = 0.9933137835985731
$ws.Range("K17").Value = 0.9937838397190122
$ws.Range("L17").Value = 0.9946018989831356
$ws.Range("M17").Value = 0.9942976182096245

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9814080470948402
$ws.Range("D18").Value = 1.02901070469991
$ws.Range("E18").Value = 0.9904784264816872
$ws.Range("F18").Value = 0.9814080470948402
$ws.Range("G18").Value = 1.005725158357515
$ws.Range("H18").Value = 1.002644322090684
$ws.Range("I18").Value = 0.9873744961291371
$ws.Range("J18").Value = 1.02901070469991
$ws.Range("K18").Value = 1.009744565590798
$ws.Range("L18").Value = 0.9955763063428194
$ws.Range("M18").Value = 0.9994401924756288

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.976643070494326
$ws.Range("D19").Value = 1.069272647050894
$ws.Range("E19").Value = 0.9815719998120205
$ws.Range("F19").Value = 0.976643070494326
$ws.Range("G19").Value = 1.029643045993584
$ws.Range("H19").Value = 0.9738047511993261
$ws.Range("I19").Value = 0.9786455506607077
$ws.Range("J19").Value = 1.069272647050894
$ws.Range("K19").Value = 1.025422323431457
$ws.Range("L19").Value = 1.001032696962892
$ws.Range("M19").Value = 1.00159684420181

# New rows 17-19 need the same bold/bordered/centered style as the other
# label (A) cells in the data range; copy formats from an existing row.
$ws.Range("A10").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0
